# Crosstab updates and other fixes.
#
# Adjusts the "exact" row heights on the two frequency tables (freq6.docx)
# and the line spacing of the trailing paragraph, matching Word's
# recalculated row metrics for this layout.
$d = $word.ActiveDocument

# Fix up the trailing paragraph's line spacing first (260 -> 264 twips),
# while Paragraphs/Tables collections are still in their pristine,
# freshly-opened state.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Format.LineSpacingRule = 5   # wdLineSpaceExactly
$p.Format.LineSpacing = 13.2    # 264 twips = 13.2 pt

# --- Table 1 ("My first Frequency Table") -------------------------------
# Both rows use a well-formed <w:trPr><w:trHeight .../></w:trPr>, so a
# plain Row.Height assignment (Word COM units are points = twips / 20)
# updates the existing element in place.
$t1 = $d.Tables.Item(1)
$t1.Rows.Item(1).HeightRule = 2  # wdRowHeightExactly
$t1.Rows.Item(1).Height = 12.8   # 252 -> 256 twips
$t1.Rows.Item(2).HeightRule = 2  # wdRowHeightExactly
$t1.Rows.Item(2).Height = 12.8   # 252 -> 256 twips

# --- Table 2 ("Eye Color" crosstab) --------------------------------------
# Several of this table's <w:trPr> elements are placed *before* their
# <w:tr> instead of inside it (a quirk of how this table was produced).
# Because of that, they are not addressable through the Rows collection's
# Height/HeightRule properties - those only ever attach a new, well-formed
# <w:trPr> onto the <w:tr>, leaving the pre-existing stray one untouched.
# Replace the whole table via InsertXML using the same markup with just
# the exact-height values updated, so the existing structure (stray
# <w:trPr> placement included) is preserved exactly:
#   - header row: 504 -> 513 twips
#   - each of the three data rows: 252 -> 256 twips
$t2 = $d.Tables.Item(2)
$xml2 = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:tbl><w:tblPr><w:tblInd w:w="2114" w:type="dxa"/><w:tblStyle w:val="TableGrid"/><w:tblW w:w="5133" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblCellMar>
                           <w:left w:w="32" w:type="dxa"/>
                           <w:right w:w="32" w:type="dxa"/>
                           </w:tblCellMar></w:tblPr>
<w:tblGrid>
<w:gridCol w:w="749"/>
<w:gridCol w:w="1137"/>
<w:gridCol w:w="878"/>
<w:gridCol w:w="1184"/>
<w:gridCol w:w="1184"/>
</w:tblGrid>
<w:tr>
<w:tc><w:tcPr><w:gridSpan w:val="5"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders>
<w:vAlign w:val="bottom"/><w:tcW w:w="5133" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Eye Color</w:t></w:r></w:p>
</w:tc>
</w:tr>
<w:tr><w:trPr><w:trHeight w:hRule="exact" w:val="513"/></w:trPr>
<w:tc><w:tcPr><w:tcW w:w="749"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:jc w:val="left"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Eyes</w:t></w:r></w:p>
</w:tc>
<w:tc><w:tcPr><w:tcW w:w="1137"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Frequency</w:t></w:r></w:p>
</w:tc>
<w:tc><w:tcPr><w:tcW w:w="878"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Percent</w:t></w:r></w:p>
</w:tc>
<w:tc><w:tcPr><w:tcW w:w="1184"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Cumulative</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Frequency</w:t></w:r></w:p>
</w:tc>
<w:tc><w:tcPr><w:tcW w:w="1184"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Cumulative</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Percent</w:t></w:r></w:p>
</w:tc>
</w:tr>

<w:trPr><w:trHeight w:hRule="exact" w:val="256"/></w:trPr><w:tr><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="left"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">blue</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">222</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">29.13</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">222</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">29.13</w:t></w:r></w:p>
</w:tc></w:tr>
<w:trPr><w:trHeight w:hRule="exact" w:val="256"/></w:trPr><w:tr><w:tc><w:p><w:pPr><w:jc w:val="left"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">brown</w:t></w:r></w:p>
</w:tc><w:tc><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">341</w:t></w:r></w:p>
</w:tc><w:tc><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">44.75</w:t></w:r></w:p>
</w:tc><w:tc><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">563</w:t></w:r></w:p>
</w:tc><w:tc><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">73.88</w:t></w:r></w:p>
</w:tc></w:tr>
<w:trPr><w:trHeight w:hRule="exact" w:val="256"/></w:trPr><w:tr><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="left"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">green</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">199</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">26.12</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">762</w:t></w:r></w:p>
</w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">100.00</w:t></w:r></w:p>
</w:tc></w:tr>
</w:tbl>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$t2.Range.InsertXML($xml2)
